$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

$ws.Range("A7").Value = "TimeSecToGetOneEnergy"
$ws.Range("B7").Value = 576

$ws.Range("A8").Value = "RequiredEnergyToPlay"
$ws.Range("B8").Value = 5

$ws.Range("A9").Value = "RefillEnergyDiamond"
$ws.Range("B9").Value = 30

$ws.Activate()
$ws.Range("B9").Select()
